$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '59.381.55'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +2.43%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.981.12'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +1.27%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '562.41'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +1.44%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '138.13'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +3.63%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.524'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +2.25%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.973.10'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +1.07%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +3.69%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.37'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +11.97%  '
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +1.28%  '
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +3.78%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '33.71'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +2.38%  '
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.37%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.472.65'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +1.30%  '
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +1.85%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.978.04'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +1.59%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '59.402.14'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +2.54%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '437.01'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +4.70%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.56'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +1.80%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.718'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +3.33%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.57%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.02'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.30%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '79.89'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +1.15%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.06%  '
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +9.49%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +0.17%  '
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +2.07%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.75'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +3.62%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.22'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +4.54%  '
$ws.Range("B32").Value = 'Hedera'
$ws.Range("C32").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.106'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +8.29%  '
$ws.Range("B33").Value = 'EthereumClassic'
$ws.Range("C33").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '25.69'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +0.74%  '
$ws.Range("B34").Value = 'PEPE'
$ws.Range("C34").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0₃0768'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +9.26%  '
$ws.Range("B35").Value = 'Filecoin'
$ws.Range("C35").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.89'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +3.56%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.981'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.03%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '48.55'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +0.65%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -0.97%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +4.31%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '400.05'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +5.11%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0352'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +1.19%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.749.20'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +1.76%  '
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -3.21%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.251'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +5.76%  '
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.03%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '35.03'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +20.19%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '122.75'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -1.28%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +2.83%  '
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +1.90%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '23.33'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +1.73%  '
